$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has an extra, redundant column (M) whose data duplicates/overlaps
# with column N. Remove column M entirely so column N's data shifts left to
# become the new column M (matches "remove column from alcohol data").
$ws.Columns.Item(13).Delete() | Out-Null

# Excel leaves the selection sitting on the column that used to be to the
# right of the deleted one (now column M) after a column delete.
$ws.Range("M1").Select() | Out-Null
